# Update the "Skeletal: Hands/Fingers and Feet/Toes" category:
#  - rename the header to "Skeletal: Hands and Feet"
#  - remove "Abnormal finger phalanx morphology" from the (alphabetised) list,
#    which shifts every following term up by one row
#  - add "Metabolic acidosis" to the (alphabetised) "Metabolism" list, right
#    before "Metabolic ketoacidosis", which now moves down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column W: "Skeletal: Hands/Fingers and Feet/Toes" -> "Skeletal: Hands and Feet"
$ws.Range("W1").Value = "Skeletal: Hands and Feet"

$handsFeet = @(
    "Aplasia/Hypoplasia of metatarsal bones",
    "Aplasia/Hypoplasia of the hallux",
    "Aplasia/Hypoplasia of the radius",
    "Brachydactyly",
    "Broad hallux",
    "Broad thumb",
    "Carpal bone hypoplasia",
    "Clubbing of fingers",
    "Clubbing of toes",
    "Drumstick terminal phalanges",
    "Foot monodactyly",
    "Hammertoe",
    "Hand monodactyly",
    "Large carpal bones",
    "Large hands",
    "Metatarsus adductus",
    "Patellar dislocation",
    "Pes planus",
    "Polydactyly",
    "Postaxial hand polydactyly",
    "Prominent fingertip pads",
    "Proximal placement of thumb",
    "Short distal phalanx of finger",
    "Short metacarpal",
    "Short metatarsal",
    "Short phalanx of finger",
    "Split foot",
    "Split hand",
    "Syndactyly",
    "Synostosis of metacarpals/metatarsals",
    "Talipes equinovalgus",
    "Talipes equinovarus",
    "Tapered finger",
    "Triphalangeal thumb"
)

for ($i = 0; $i -lt $handsFeet.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 23).Value = $handsFeet[$i]
}

# The list is now one entry shorter, so the last (previously occupied) row is cleared
$ws.Cells.Item(36, 23).Value = ""

# --- Column AA: "Metabolism" list gains "Metabolic acidosis" before "Metabolic ketoacidosis"
$ws.Range("AA33").Value = "Metabolic acidosis"
$ws.Range("AA34").Value = "Metabolic ketoacidosis"
